# sample3.xlsx edit: reshuffle weekday schedule cells, drop the extra
# "AXS"/VLK legend row, and replace its link with a fresh Google Meet URL
# used by CAOS LAB. ("Pushing before removing the buttons-link")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Monday (row 2) ---
$ws.Range("B2").Value = "SSK"
$ws.Range("E2").Value = ""

# --- Tuesday (row 3) ---
$ws.Range("B3").Value = "DDB"
$ws.Range("F3").Value = "AAB DS"

# --- Wednesday (row 4) ---
$ws.Range("B4").Value = "VLK"
$ws.Range("C4").Value = "AAB DS"
$ws.Range("D4").Value = "AAB DS LAB"
$ws.Range("E4").Value = "CAOS LAB"
$ws.Range("G4").Value = "VDP DCAN LAB"
$ws.Range("H4").Value = "BCD"
$ws.Range("I4").Value = "SSK"

# --- Thursday (row 5) ---
$ws.Range("B5").Value = "VDP DCAN LAB"
$ws.Range("C5").Value = "VLK"
$ws.Range("E5").Value = "CAOS"
$ws.Range("F5").Value = ""
$ws.Range("H5").Value = ""
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = "CAOS LAB"

# --- Friday (row 6) ---
$ws.Range("B6").Value = "AAB DS"
$ws.Range("C6").Value = "CAOS"
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = "SSK"
$ws.Range("F6").Value = "DDB"
$ws.Range("H6").Value = "BCD"
$ws.Range("I6").Value = ""
$ws.Range("J6").Value = ""

# --- Legend / links table (column Q/R) ---
# CAOS LAB's link is replaced with the new meeting URL.
$ws.Range("R9").Value = "https://meet.google.com/khj-syqw-kvz"

# The old row 11 (AXS / duplicate link) is no longer needed now that its
# link moved up to row 9 - remove the whole row.
$ws.Rows("11:11").Delete()

# Column widths: best-fit on the code column (B) and the legend-code
# column (Q) so the new longer labels aren't clipped.
$ws.Columns("B:B").ColumnWidth = 13.6640625
$ws.Columns("Q:Q").ColumnWidth = 13.6640625

# Restore the selection to the cell the author left active.
$ws.Range("F3").Select()
